# Inserts a new weekly record for "Apio" (Americana, Primera) at row 39,
# pushing the existing rows 39:172 down to 40:173.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(39).Insert()

$ws.Range("A39").Value = 5
$ws.Range("B39").Value = "Macroferia Regional de Talca"
$ws.Range("C39").Value = "Maule"
$ws.Range("D39").Value = 44701
$ws.Range("E39").Value = 7
$ws.Range("F39").Value = 100112017
$ws.Range("G39").Value = "Apio"
$ws.Range("H39").Value = "Americana (o)"
$ws.Range("I39").Value = "Primera"
$ws.Range("J39").Value = 700
$ws.Range("K39").Value = 7000
$ws.Range("L39").Value = 7000
$ws.Range("M39").Value = 7000
$ws.Range("N39").Value = "`$/docena de matas"
$ws.Range("O39").Value = "Provincia del Elquí"
$ws.Range("P39").Value = 1167
$ws.Range("Q39").Value = 6
$ws.Range("R39").Value = "Hortaliza"
